$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 3,4,6
foreach ($r in $rows) {
    $ws.Range("A$r").Value = "Rajesh"
    $ws.Range("B$r").Value = 35920
    $ws.Range("C$r").Value = "Y"
    $ws.Range("D$r").Value = 1000.1
    $ws.Range("E$r").Value = 100000
}

$ws.Range("B5").Value = 35920
$ws.Range("C5").Value = "Y"
$ws.Range("D5").Value = 1000.1
$ws.Range("E5").Value = 100000

$ws.Range("B2").Copy()
$ws.Range("B3:B6").PasteSpecial(-4122)

$ws.Range("K13").Select()
